$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 159, pushing all existing rows (159-243) down
# to become rows 161-245.
$ws.Rows("159:160").Insert()

# Fill in the new row 159 with the new weekly data point ("1a (guarda)")
$ws.Range("A159").Value = 11
$ws.Range("B159").Value = "Vega Monumental Concepción"
$ws.Range("C159").Value = "Bíobío"
$ws.Range("D159").Value = (Get-Date -Year 2021 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E159").Value = 8
$ws.Range("F159").Value = 100112004
$ws.Range("G159").Value = "Cebolla"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "1a (guarda)"
$ws.Range("J159").Value = 1000
$ws.Range("K159").Value = 6000
$ws.Range("L159").Value = 6500
$ws.Range("M159").Value = 6250
$ws.Range("N159").Value = "`$/malla 18 kilos"
$ws.Range("O159").Value = "Región Metropolitana"
$ws.Range("P159").Value = 347
$ws.Range("Q159").Value = 18
$ws.Range("R159").Value = "Hortaliza"

# Fill in the new row 160 with the new weekly data point ("2a (guarda)")
$ws.Range("A160").Value = 11
$ws.Range("B160").Value = "Vega Monumental Concepción"
$ws.Range("C160").Value = "Bíobío"
$ws.Range("D160").Value = (Get-Date -Year 2021 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E160").Value = 8
$ws.Range("F160").Value = 100112004
$ws.Range("G160").Value = "Cebolla"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "2a (guarda)"
$ws.Range("J160").Value = 500
$ws.Range("K160").Value = 5500
$ws.Range("L160").Value = 5500
$ws.Range("M160").Value = 5500
$ws.Range("N160").Value = "`$/malla 18 kilos"
$ws.Range("O160").Value = "Región Metropolitana"
$ws.Range("P160").Value = 306
$ws.Range("Q160").Value = 18
$ws.Range("R160").Value = "Hortaliza"
